# Apply scheduled market-data refresh to the Leve profit sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Each sheet is an Excel Table (Table_<code>) spanning A1:N141 with columns:
#   H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#   K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
# Updates below mirror refreshed market-board prices and their dependent profit calcs.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 9832.5
$ws.Range("I8").Value = 9799.200000000001
$ws.Range("K8").Value = 29397.6
$ws.Range("M8").Value = -29258.6
$ws.Range("H62").Value = 60318984
$ws.Range("I62").Value = 66667932
$ws.Range("K62").Value = 66667932
$ws.Range("M62").Value = -66667308
$ws.Range("H65").Value = 60318984
$ws.Range("I65").Value = 66667932
$ws.Range("K65").Value = 333339660
$ws.Range("M65").Value = -333336540
$ws.Range("H70").Value = 4172.353
$ws.Range("I70").Value = 3699.25
$ws.Range("J70").Value = 4317.923
$ws.Range("K70").Value = 11097.75
$ws.Range("L70").Value = 12953.769
$ws.Range("M70").Value = -10827.75
$ws.Range("N70").Value = -13493.769
$ws.Range("H73").Value = 4172.353
$ws.Range("I73").Value = 3699.25
$ws.Range("J73").Value = 4317.923
$ws.Range("K73").Value = 11097.75
$ws.Range("L73").Value = 12953.769
$ws.Range("M73").Value = -10161.75
$ws.Range("N73").Value = -14825.769
$ws.Range("H74").Value = 8242.083000000001
$ws.Range("I74").Value = 6581.2
$ws.Range("J74").Value = 9428.429
$ws.Range("K74").Value = 6581.2
$ws.Range("L74").Value = 9428.429
$ws.Range("M74").Value = -5645.2
$ws.Range("N74").Value = -11300.429
$ws.Range("H77").Value = 8242.083000000001
$ws.Range("I77").Value = 6581.2
$ws.Range("J77").Value = 9428.429
$ws.Range("K77").Value = 32906
$ws.Range("L77").Value = 47142.145
$ws.Range("M77").Value = -28226
$ws.Range("N77").Value = -56502.145
$ws.Range("H111").Value = 2754.0715
$ws.Range("I111").Value = 2321.625
$ws.Range("K111").Value = 6964.875
$ws.Range("M111").Value = -3897.875
$ws.Range("H141").Value = 964.6667
$ws.Range("I141").Value = 964.6667
$ws.Range("K141").Value = 2894.0001
$ws.Range("M141").Value = 2285.9999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 782.25
$ws.Range("I88").Value = 712
$ws.Range("J88").Value = 832.4286
$ws.Range("K88").Value = 712
$ws.Range("L88").Value = 832.4286
$ws.Range("M88").Value = -306
$ws.Range("N88").Value = -1644.4286
$ws.Range("H91").Value = 782.25
$ws.Range("I91").Value = 712
$ws.Range("J91").Value = 832.4286
$ws.Range("K91").Value = 712
$ws.Range("L91").Value = 832.4286
$ws.Range("M91").Value = 692
$ws.Range("N91").Value = -3640.4286
$ws.Range("H132").Value = 4692.467
$ws.Range("I132").Value = 3068.2307
$ws.Range("K132").Value = 9204.6921
$ws.Range("M132").Value = -6674.6921
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 47621076
$ws.Range("I105").Value = 62501876
$ws.Range("J105").Value = 2524.4
$ws.Range("K105").Value = 62501876
$ws.Range("L105").Value = 2524.4
$ws.Range("M105").Value = -62500129
$ws.Range("N105").Value = -6018.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28575594
$ws.Range("H34").Value = 28575594
$ws.Range("H58").Value = 6194.3213
$ws.Range("I58").Value = 3883.5
$ws.Range("K58").Value = 3883.5
$ws.Range("M58").Value = -3680.5
$ws.Range("H99").Value = 1933.3334
$ws.Range("I99").Value = 1933.3334
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1933.3334
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -435.3334
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 1933.3334
$ws.Range("I126").Value = 1933.3334
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5800.0002
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3330.0002
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 211783
$ws.Range("I132").Value = 11974.5
$ws.Range("J132").Value = 311687.25
$ws.Range("K132").Value = 35923.5
$ws.Range("L132").Value = 935061.75
$ws.Range("M132").Value = -33393.5
$ws.Range("N132").Value = -940121.75
$ws.Range("H136").Value = 6194.3213
$ws.Range("I136").Value = 3883.5
$ws.Range("K136").Value = 11650.5
$ws.Range("M136").Value = -9100.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4586.25
$ws.Range("I3").Value = 4586.25
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 13758.75
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -13646.75
$ws.Range("N3").ClearContents()
$ws.Range("H5").Value = 947.1
$ws.Range("I5").Value = 598.6667
$ws.Range("J5").Value = 1096.4286
$ws.Range("K5").Value = 1796.0001
$ws.Range("L5").Value = 3289.2858
$ws.Range("M5").Value = -1684.0001
$ws.Range("N5").Value = -3513.2858
$ws.Range("H44").Value = 2302.1
$ws.Range("I44").Value = 586.8333
$ws.Range("K44").Value = 1760.4999
$ws.Range("M44").Value = -1362.4999
$ws.Range("H74").Value = 11631.5
$ws.Range("H77").Value = 11631.5
$ws.Range("H129").Value = 23812716
$ws.Range("J129").Value = 166667460
$ws.Range("L129").Value = 500002380
$ws.Range("N129").Value = -500012380
$ws.Range("H131").Value = 7580478.5
$ws.Range("J131").Value = 5497.5557
$ws.Range("L131").Value = 16492.6671
$ws.Range("N131").Value = -26572.6671
$ws.Range("H135").Value = 947.1
$ws.Range("I135").Value = 598.6667
$ws.Range("J135").Value = 1096.4286
$ws.Range("K135").Value = 5388.0003
$ws.Range("L135").Value = 9867.857399999999
$ws.Range("M135").Value = -2853.0003
$ws.Range("N135").Value = -14937.8574
$ws.Range("H137").Value = 2731.5
$ws.Range("J137").Value = 3982.3333
$ws.Range("L137").Value = 11946.9999
$ws.Range("N137").Value = -22146.9999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4978.75
$ws.Range("I132").Value = 2071.6667
$ws.Range("K132").Value = 6215.000100000001
$ws.Range("M132").Value = -3685.000100000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 8465.666999999999
$ws.Range("I122").Value = 8465.666999999999
$ws.Range("K122").Value = 25397.001
$ws.Range("M122").Value = -22947.001
$ws.Range("H132").Value = 4989.4546
$ws.Range("I132").Value = 3774.2856
$ws.Range("K132").Value = 11322.8568
$ws.Range("M132").Value = -8792.856800000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 14737.321
$ws.Range("I81").Value = 2583.3333
$ws.Range("J81").Value = 16195.8
$ws.Range("K81").Value = 5166.6666
$ws.Range("L81").Value = 32391.6
$ws.Range("M81").Value = -4105.6666
$ws.Range("N81").Value = -34513.6
$ws.Range("H84").Value = 14737.321
$ws.Range("I84").Value = 2583.3333
$ws.Range("J84").Value = 16195.8
$ws.Range("K84").Value = 25833.333
$ws.Range("L84").Value = 161958
$ws.Range("M84").Value = -20529.333
$ws.Range("N84").Value = -172566
$ws.Range("H86").Value = 65792.336
$ws.Range("J86").Value = 65792.336
$ws.Range("L86").Value = 65792.336
$ws.Range("N86").Value = -68038.336
$ws.Range("H89").Value = 65792.336
$ws.Range("J89").Value = 65792.336
$ws.Range("L89").Value = 328961.68
$ws.Range("N89").Value = -340193.68
$ws.Range("H100").Value = 2361.6177
$ws.Range("I100").Value = 2606.6785
$ws.Range("J100").Value = 1218
$ws.Range("K100").Value = 5213.357
$ws.Range("L100").Value = 2436
$ws.Range("M100").Value = -4672.357
$ws.Range("N100").Value = -3518
$ws.Range("H113").Value = 518.15
$ws.Range("I113").Value = 404.9375
$ws.Range("J113").Value = 971
$ws.Range("K113").Value = 1214.8125
$ws.Range("L113").Value = 2913
$ws.Range("M113").Value = 955.1875
$ws.Range("N113").Value = -7253
$ws.Range("H132").Value = 9776.846
$ws.Range("I132").Value = 8000
$ws.Range("K132").Value = 24000
$ws.Range("M132").Value = -21470
